$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge-weight table following Dr Hou's advice:
# replaces existing rows 2-13 values and appends new rows 14-17
# (sCs -> {ECs,FAPs,M2,sCs} target-cluster combinations).
$data = New-Object 'object[,]' 16,20
$data[0,0] = 'ECs'
$data[0,1] = 'Sertad1'
$data[0,2] = 'Ar'
$data[0,3] = 'ECs'
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 21.510992
$data[0,7] = 64.53297600000001
$data[0,8] = 0.3874081946303762
$data[0,9] = 0.3874081946303762
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.8164263333333333
$data[0,13] = 2.449279
$data[0,14] = 0.05945196387955778
$data[0,15] = 0.05945196387955777
$data[0,16] = 17.56214032492267
$data[0,17] = 158.059262924304
$data[0,18] = 0.02303217799380982
$data[0,19] = 0.02303217799380981
$data[1,0] = 'ECs'
$data[1,1] = 'Sertad1'
$data[1,2] = 'Ar'
$data[1,3] = 'FAPs'
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 21.510992
$data[1,7] = 64.53297600000001
$data[1,8] = 0.3874081946303762
$data[1,9] = 0.3874081946303762
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 10.32138733333333
$data[1,13] = 30.964162
$data[1,14] = 0.7516008755167443
$data[1,15] = 0.7516008755167441
$data[1,16] = 222.0232803562347
$data[1,17] = 1998.209523206112
$data[1,18] = 0.291176338266552
$data[1,19] = 0.291176338266552
$data[2,0] = 'ECs'
$data[2,1] = 'Sertad1'
$data[2,2] = 'Ar'
$data[2,3] = 'M2'
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 21.510992
$data[2,7] = 64.53297600000001
$data[2,8] = 0.3874081946303762
$data[2,9] = 0.3874081946303762
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.037004
$data[2,13] = 0.111012
$data[2,14] = 0.002694622137452478
$data[2,15] = 0.002694622137452478
$data[2,16] = 0.7959927479680001
$data[2,17] = 7.163934731712001
$data[2,18] = 0.00104391869748151
$data[2,19] = 0.00104391869748151
$data[3,0] = 'ECs'
$data[3,1] = 'Sertad1'
$data[3,2] = 'Ar'
$data[3,3] = 'sCs'
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 21.510992
$data[3,7] = 64.53297600000001
$data[3,8] = 0.3874081946303762
$data[3,9] = 0.3874081946303762
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 2.55772
$data[3,13] = 7.67316
$data[3,14] = 0.1862525384662456
$data[3,15] = 0.1862525384662456
$data[3,16] = 55.01909445824001
$data[3,17] = 495.1718501241601
$data[3,18] = 0.07215575967253293
$data[3,19] = 0.0721557596725329
$data[4,0] = 'FAPs'
$data[4,1] = 'Sertad1'
$data[4,2] = 'Ar'
$data[4,3] = 'ECs'
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 15.69618166666667
$data[4,7] = 47.088545
$data[4,8] = 0.2826847502929545
$data[4,9] = 0.2826847502929545
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.8164263333333333
$data[4,13] = 2.449279
$data[4,14] = 0.05945196387955778
$data[4,15] = 0.05945196387955777
$data[4,16] = 12.81477604545055
$data[4,17] = 115.332984409055
$data[4,18] = 0.01680616356371854
$data[4,19] = 0.01680616356371854
$data[5,0] = 'FAPs'
$data[5,1] = 'Sertad1'
$data[5,2] = 'Ar'
$data[5,3] = 'FAPs'
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 15.69618166666667
$data[5,7] = 47.088545
$data[5,8] = 0.2826847502929545
$data[5,9] = 0.2826847502929545
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 10.32138733333333
$data[5,13] = 30.964162
$data[5,14] = 0.7516008755167443
$data[5,15] = 0.7516008755167441
$data[5,16] = 162.0063706360322
$data[5,17] = 1458.05733572429
$data[5,18] = 0.2124661058154169
$data[5,19] = 0.2124661058154168
$data[6,0] = 'FAPs'
$data[6,1] = 'Sertad1'
$data[6,2] = 'Ar'
$data[6,3] = 'M2'
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 15.69618166666667
$data[6,7] = 47.088545
$data[6,8] = 0.2826847502929545
$data[6,9] = 0.2826847502929545
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.037004
$data[6,13] = 0.111012
$data[6,14] = 0.002694622137452478
$data[6,15] = 0.002694622137452478
$data[6,16] = 0.5808215063933333
$data[6,17] = 5.227393557539999
$data[6,18] = 0.0007617285860596213
$data[6,19] = 0.0007617285860596211
$data[7,0] = 'FAPs'
$data[7,1] = 'Sertad1'
$data[7,2] = 'Ar'
$data[7,3] = 'sCs'
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 15.69618166666667
$data[7,7] = 47.088545
$data[7,8] = 0.2826847502929545
$data[7,9] = 0.2826847502929545
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 2.55772
$data[7,13] = 7.67316
$data[7,14] = 0.1862525384662456
$data[7,15] = 0.1862525384662456
$data[7,16] = 40.14643777246667
$data[7,17] = 361.3179399522
$data[7,18] = 0.05265075232775956
$data[7,19] = 0.05265075232775954
$data[8,0] = 'M2'
$data[8,1] = 'Sertad1'
$data[8,2] = 'Ar'
$data[8,3] = 'ECs'
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 12.36237333333333
$data[8,7] = 37.08712
$data[8,8] = 0.2226436016718045
$data[8,9] = 0.2226436016718045
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.8164263333333333
$data[8,13] = 2.449279
$data[8,14] = 0.05945196387955778
$data[8,15] = 0.05945196387955777
$data[8,16] = 10.09296713183111
$data[8,17] = 90.83670418647999
$data[8,18] = 0.01323659936460677
$data[8,19] = 0.01323659936460677
$data[9,0] = 'M2'
$data[9,1] = 'Sertad1'
$data[9,2] = 'Ar'
$data[9,3] = 'FAPs'
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 12.36237333333333
$data[9,7] = 37.08712
$data[9,8] = 0.2226436016718045
$data[9,9] = 0.2226436016718045
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 10.32138733333333
$data[9,13] = 30.964162
$data[9,14] = 0.7516008755167443
$data[9,15] = 0.7516008755167441
$data[9,16] = 127.5968435326044
$data[9,17] = 1148.37159179344
$data[9,18] = 0.1673391259447295
$data[9,19] = 0.1673391259447295
$data[10,0] = 'M2'
$data[10,1] = 'Sertad1'
$data[10,2] = 'Ar'
$data[10,3] = 'M2'
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 12.36237333333333
$data[10,7] = 37.08712
$data[10,8] = 0.2226436016718045
$data[10,9] = 0.2226436016718045
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.037004
$data[10,13] = 0.111012
$data[10,14] = 0.002694622137452478
$data[10,15] = 0.002694622137452478
$data[10,16] = 0.4574572628266667
$data[10,17] = 4.11711536544
$data[10,18] = 0.000599940377826996
$data[10,19] = 0.0005999403778269959
$data[11,0] = 'M2'
$data[11,1] = 'Sertad1'
$data[11,2] = 'Ar'
$data[11,3] = 'sCs'
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 12.36237333333333
$data[11,7] = 37.08712
$data[11,8] = 0.2226436016718045
$data[11,9] = 0.2226436016718045
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 2.55772
$data[11,13] = 7.67316
$data[11,14] = 0.1862525384662456
$data[11,15] = 0.1862525384662456
$data[11,16] = 31.61948952213333
$data[11,17] = 284.5754056992
$data[11,18] = 0.04146793598464123
$data[11,19] = 0.04146793598464122
$data[12,0] = 'sCs'
$data[12,1] = 'Sertad1'
$data[12,2] = 'Ar'
$data[12,3] = 'ECs'
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 5.955845333333333
$data[12,7] = 17.867536
$data[12,8] = 0.1072634534048647
$data[12,9] = 0.1072634534048647
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 0.8164263333333333
$data[12,13] = 2.449279
$data[12,14] = 0.05945196387955778
$data[12,15] = 0.05945196387955777
$data[12,16] = 4.862508967393778
$data[12,17] = 43.762580706544
$data[12,18] = 0.006377022957422647
$data[12,19] = 0.006377022957422646
$data[13,0] = 'sCs'
$data[13,1] = 'Sertad1'
$data[13,2] = 'Ar'
$data[13,3] = 'FAPs'
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 5.955845333333333
$data[13,7] = 17.867536
$data[13,8] = 0.1072634534048647
$data[13,9] = 0.1072634534048647
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 10.32138733333333
$data[13,13] = 30.964162
$data[13,14] = 0.7516008755167443
$data[13,15] = 0.7516008755167441
$data[13,16] = 61.47258658275911
$data[13,17] = 553.253279244832
$data[13,18] = 0.08061930549004584
$data[13,19] = 0.08061930549004583
$data[14,0] = 'sCs'
$data[14,1] = 'Sertad1'
$data[14,2] = 'Ar'
$data[14,3] = 'M2'
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 5.955845333333333
$data[14,7] = 17.867536
$data[14,8] = 0.1072634534048647
$data[14,9] = 0.1072634534048647
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.037004
$data[14,13] = 0.111012
$data[14,14] = 0.002694622137452478
$data[14,15] = 0.002694622137452478
$data[14,16] = 0.2203901007146667
$data[14,17] = 1.983510906432
$data[14,18] = 0.0002890344760843509
$data[14,19] = 0.0002890344760843509
$data[15,0] = 'sCs'
$data[15,1] = 'Sertad1'
$data[15,2] = 'Ar'
$data[15,3] = 'sCs'
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 5.955845333333333
$data[15,7] = 17.867536
$data[15,8] = 0.1072634534048647
$data[15,9] = 0.1072634534048647
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 2.55772
$data[15,13] = 7.67316
$data[15,14] = 0.1862525384662456
$data[15,15] = 0.1862525384662456
$data[15,16] = 15.23338472597334
$data[15,17] = 137.10046253376
$data[15,18] = 0.01997809048131192
$data[15,19] = 0.01997809048131191

$ws.Range("A2:T17").Value = $data
